# Updated Logic Command Class Diagram
#
# Resizes/repositions several shapes in the Logic component diagram and
# renames the Add/Clear/Find command boxes to their "...Person..." variants.
#
# Shape.Left/Top/Width/Height on this COM surface are expressed in points
# (1 pt = 12700 EMU); we add a half-EMU (in point units) before assigning so
# that the internal points->EMU conversion (which truncates) lands on the
# exact target EMU value instead of one EMU short.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function EmuToPt($emu) {
    return ($emu / 12700.0) + (0.5 / 12700.0)
}

function SetBounds($shape, $left, $top, $width, $height) {
    if ($null -ne $left)   { $shape.Left   = EmuToPt $left }
    if ($null -ne $top)    { $shape.Top    = EmuToPt $top }
    if ($null -ne $width)  { $shape.Width  = EmuToPt $width }
    if ($null -ne $height) { $shape.Height = EmuToPt $height }
}

# --- id=2 "Logic" rounded rectangle: widen / shift left -------------------
$logicBox = $s.Shapes.Item(2)
SetBounds $logicBox 2557294 $null 6358106 $null

# --- id=107 "AddCommand" -> "AddPersonCommand" -----------------------------
$addCmd = $s.Shapes.Item(8)
SetBounds $addCmd 7238999 $null 1447798 $null
$addCmd.TextFrame.TextRange.Text = "AddPersonCommand"

# --- id=113 "ClearCommand" -> "ClearPersonCommand" -------------------------
$clearCmd = $s.Shapes.Item(9)
SetBounds $clearCmd $null $null 1447800 $null
$clearCmd.TextFrame.TextRange.Text = "ClearPersonCommand"

# --- id=128 "..." command box: shift / widen (no rename) -------------------
$ellipsisCmd1 = $s.Shapes.Item(10)
SetBounds $ellipsisCmd1 7238999 $null 1447797 $null

# --- id=129 connector: tiny width correction -------------------------------
$conn129 = $s.Shapes.Item(11)
SetBounds $conn129 $null $null 1056756 $null

# --- id=134 connector: tiny width correction -------------------------------
$conn134 = $s.Shapes.Item(13)
SetBounds $conn134 $null $null 1056756 $null

# --- id=144 "FindCommand" -> "FindPersonCommand" ----------------------------
$findCmd = $s.Shapes.Item(15)
SetBounds $findCmd 3049433 $null 1455096 $null
$findCmd.TextFrame.TextRange.Text = "FindPersonCommand"

# --- id=145 "UndoCommand": shift / widen (no rename) ------------------------
$undoCmd = $s.Shapes.Item(16)
SetBounds $undoCmd 3049433 $null 1455096 $null

# --- id=167 "..." command box: shift / widen (no rename) -------------------
$ellipsisCmd2 = $s.Shapes.Item(18)
SetBounds $ellipsisCmd2 3048128 $null 1455096 $null

# --- id=30 connector: shrink width/height -----------------------------------
$conn30 = $s.Shapes.Item(25)
SetBounds $conn30 $null $null 413532 4120
